$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.812.26"
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.313.55"
$ws.Range("E3").Value = "  +0.96%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.58"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.58"
$ws.Range("E6").Value = "  +1.81%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.631"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +1.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.27"
$ws.Range("E10").Value = "  +2.99%  "
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.58"
$ws.Range("E12").Value = "  +3.85%  "
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.977"
$ws.Range("E14").Value = "  +1.65%  "
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.660.05"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.313.45"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.739.85"
$ws.Range("E18").Value = "  +1.67%  "
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.26"
$ws.Range("E21").Value = "  +31.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.89"
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("E23").Value = "  -3.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.30"
$ws.Range("E24").Value = "  -4.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.27"
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.92"
$ws.Range("E27").Value = "  +1.29%  "
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.71"
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.08"
$ws.Range("E30").Value = "  +6.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.27"
$ws.Range("E31").Value = "  +7.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.91"
$ws.Range("E32").Value = "  +1.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0891"
$ws.Range("E33").Value = "  +2.15%  "
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.58"
$ws.Range("E35").Value = "  -9.39%  "
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("E37").Value = "  +2.15%  "
$ws.Range("E38").Value = "  +1.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.72"
$ws.Range("E39").Value = "  +1.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.74"
$ws.Range("E40").Value = "  -4.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.59"
$ws.Range("E41").Value = "  +9.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "98.67"
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.45"
$ws.Range("E43").Value = "  +1.47%  "
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.45"
$ws.Range("E46").Value = "  +4.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "82.06"
$ws.Range("E47").Value = "  +7.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "114.77"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.32"
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("E50").Value = "  -0.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.615.96"
$ws.Range("E51").Value = "  +5.36%  "
